$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the full Year/Payment table (1996-2018), preserving the
# original 2009-2014 payment values and filling the rest with 0.
$years = 1996..2018
$payments = @{
    2009 = 200
    2010 = 150
    2011 = 100
    2012 = 100
    2013 = 100
    2014 = 100
}

$row = 2
foreach ($year in $years) {
    $amount = 0
    if ($payments.ContainsKey($year)) {
        $amount = $payments[$year]
    }
    $ws.Cells.Item($row, 1).Value = $year
    $ws.Cells.Item($row, 2).Value = $amount
    $row = $row + 1
}

# Update the selected/visible range shown in the saved view (scrolled
# so row 7 is at the top, with E7:G20 selected).
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E7:G20").Select()
